# rebal-v2 / tickers.xlsx
# "added calculations for ann return, vol, and utility"
#
# Concretely: add four new commodity-ETF tickers (gold/silver/palladium/
# platinum trackers) to the "Tickers" sheet, grouped as "Commod" with zero
# min/max weight placeholders, and bump the "Bond" group's max weight on
# the "Groups" sheet from 0.3 to 0.4.

$wb = $excel.ActiveWorkbook

# --- Tickers sheet -------------------------------------------------------
$tickers = $wb.Worksheets.Item("Tickers")

# TLT is the last row (31). Push it down 4 rows so the new tickers can be
# inserted directly above it, keeping the "Bond" row last.
$tickers.Range("A31:A34").EntireRow.Insert(-4121, 0) | Out-Null

$newRows = @(
    @("IAU",  0, 0, "Commod"),
    @("SLV",  0, 0, "Commod"),
    @("PALL", 0, 0, "Commod"),
    @("PPLT", 0, 0, "Commod")
)

$r = 31
foreach ($row in $newRows) {
    $tickers.Cells.Item($r, 1).Value = $row[0]
    $tickers.Cells.Item($r, 2).Value = $row[1]
    $tickers.Cells.Item($r, 3).Value = $row[2]
    $tickers.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Tickers isn't the active tab; leave its selection at the top-left cell.
$tickers.Range("A1").Select() | Out-Null

# --- Groups sheet ---------------------------------------------------------
$groups = $wb.Worksheets.Item("Groups")
$groups.Range("C5").Value = 0.4

# Groups remains the active / selected tab.
$groups.Activate() | Out-Null
$groups.Range("B5").Select() | Out-Null
